# Atualização de bases das ligas, do dia: 16-06-2024 às 07:16
#
# The source feed had a handful of match rows whose records (columns B:AD —
# everything except the running index in column A) were written to the wrong
# row. This script fixes the ordering by swapping/rotating the B:AD payload
# between the affected rows, leaving column A (the row index) untouched.
#
# NOTE: multi-cell Range.Value reads are not reliable in this host, so we
# read/write cell-by-cell via Cells.Item(row, col).Value2 / .Value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 30  # column AD

function Get-RowPayload($ws, [int]$row) {
    $vals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals += , ($ws.Cells.Item($row, $c).Value2)
    }
    return $vals
}

function Set-RowPayload($ws, [int]$row, $vals) {
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $c = $firstCol + $i
        $ws.Cells.Item($row, $c).Value = $vals[$i]
    }
}

function Swap-RowPayload($ws, [int]$rowA, [int]$rowB) {
    # Swap everything in columns B..AD between two rows, keeping column A put.
    $valsA = Get-RowPayload $ws $rowA
    $valsB = Get-RowPayload $ws $rowB

    Set-RowPayload $ws $rowA $valsB
    Set-RowPayload $ws $rowB $valsA
}

# Simple pairwise swaps
Swap-RowPayload $ws 175 176
Swap-RowPayload $ws 181 182
Swap-RowPayload $ws 183 184
Swap-RowPayload $ws 185 186
Swap-RowPayload $ws 187 188

# Three-way rotation: new(313) = old(315); new(314) = old(313); new(315) = old(314)
$v313 = Get-RowPayload $ws 313
$v314 = Get-RowPayload $ws 314
$v315 = Get-RowPayload $ws 315

Set-RowPayload $ws 313 $v315
Set-RowPayload $ws 314 $v313
Set-RowPayload $ws 315 $v314
